$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Vwf"
$ws.Range("C2").Value = "Tnfrsf11b"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 73.450124
$ws.Range("H2").Value = 220.350372
$ws.Range("I2").Value = 0.9139420548774678
$ws.Range("J2").Value = 0.9139420548774677
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.08241233333333334
$ws.Range("N2").Value = 0.247237
$ws.Range("O2").Value = 0.04727005612861496
$ws.Range("P2").Value = 0.04727005612861496
$ws.Range("Q2").Value = 6.053196102462667
$ws.Range("R2").Value = 54.478764922164
$ws.Range("S2").Value = 0.0432020922323596
$ws.Range("T2").Value = 0.04320209223235959

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Vwf"
$ws.Range("C3").Value = "Tnfrsf11b"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 73.450124
$ws.Range("H3").Value = 220.350372
$ws.Range("I3").Value = 0.9139420548774678
$ws.Range("J3").Value = 0.9139420548774677
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.661024
$ws.Range("N3").Value = 4.983072
$ws.Range("O3").Value = 0.9527299438713851
$ws.Range("P3").Value = 0.952729943871385
$ws.Range("Q3").Value = 122.002418766976
$ws.Range("R3").Value = 1098.021768902784
$ws.Range("S3").Value = 0.8707399626451083
$ws.Range("T3").Value = 0.8707399626451081

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Vwf"
$ws.Range("C4").Value = "Tnfrsf11b"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1722003333333334
$ws.Range("H4").Value = 0.5166010000000001
$ws.Range("I4").Value = 0.002142693816245315
$ws.Range("J4").Value = 0.002142693816245315
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.08241233333333334
$ws.Range("N4").Value = 0.247237
$ws.Range("O4").Value = 0.04727005612861496
$ws.Range("P4").Value = 0.04727005612861496
$ws.Range("Q4").Value = 0.01419143127077778
$ws.Range("R4").Value = 0.127722881437
$ws.Range("S4").Value = 0.0001012852569603523
$ws.Range("T4").Value = 0.0001012852569603522

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Vwf"
$ws.Range("C5").Value = "Tnfrsf11b"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.1722003333333334
$ws.Range("H5").Value = 0.5166010000000001
$ws.Range("I5").Value = 0.002142693816245315
$ws.Range("J5").Value = 0.002142693816245315
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.661024
$ws.Range("N5").Value = 4.983072
$ws.Range("O5").Value = 0.9527299438713851
$ws.Range("P5").Value = 0.952729943871385
$ws.Range("Q5").Value = 0.2860288864746667
$ws.Range("R5").Value = 2.574259978272
$ws.Range("S5").Value = 0.002041408559284963
$ws.Range("T5").Value = 0.002041408559284963

# Row 6
$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "Vwf"
$ws.Range("C6").Value = "Tnfrsf11b"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.623825
$ws.Range("H6").Value = 10.871475
$ws.Range("I6").Value = 0.04509136113938133
$ws.Range("J6").Value = 0.04509136113938133
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.08241233333333334
$ws.Range("N6").Value = 0.247237
$ws.Range("O6").Value = 0.04727005612861496
$ws.Range("P6").Value = 0.04727005612861496
$ws.Range("Q6").Value = 0.2986478738416667
$ws.Range("R6").Value = 2.687830864575
$ws.Range("S6").Value = 0.002131471171974203
$ws.Range("T6").Value = 0.002131471171974203

# Row 7
$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "Vwf"
$ws.Range("C7").Value = "Tnfrsf11b"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.623825
$ws.Range("H7").Value = 10.871475
$ws.Range("I7").Value = 0.04509136113938133
$ws.Range("J7").Value = 0.04509136113938133
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.661024
$ws.Range("N7").Value = 4.983072
$ws.Range("O7").Value = 0.9527299438713851
$ws.Range("P7").Value = 0.952729943871385
$ws.Range("Q7").Value = 6.019260296800001
$ws.Range("R7").Value = 54.1733426712
$ws.Range("S7").Value = 0.04295988996740713
$ws.Range("T7").Value = 0.04295988996740712

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Vwf"
$ws.Range("C8").Value = "Tnfrsf11b"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.871251333333333
$ws.Range("H8").Value = 8.613754
$ws.Range("I8").Value = 0.03572706485364594
$ws.Range("J8").Value = 0.03572706485364593
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.08241233333333334
$ws.Range("N8").Value = 0.247237
$ws.Range("O8").Value = 0.04727005612861496
$ws.Range("P8").Value = 0.04727005612861496
$ws.Range("Q8").Value = 0.2366265219664445
$ws.Range("R8").Value = 2.129638697698
$ws.Range("S8").Value = 0.001688820360942511
$ws.Range("T8").Value = 0.00168882036094251

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Vwf"
$ws.Range("C9").Value = "Tnfrsf11b"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.871251333333333
$ws.Range("H9").Value = 8.613754
$ws.Range("I9").Value = 0.03572706485364594
$ws.Range("J9").Value = 0.03572706485364593
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.661024
$ws.Range("N9").Value = 4.983072
$ws.Range("O9").Value = 0.9527299438713851
$ws.Range("P9").Value = 0.952729943871385
$ws.Range("Q9").Value = 4.769217374698667
$ws.Range("R9").Value = 42.922956372288
$ws.Range("S9").Value = 0.03403824449270344
$ws.Range("T9").Value = 0.03403824449270342

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Vwf"
$ws.Range("C10").Value = "Tnfrsf11b"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.2488803333333334
$ws.Range("H10").Value = 0.7466410000000001
$ws.Range("I10").Value = 0.003096825313259592
$ws.Range("J10").Value = 0.003096825313259592
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.08241233333333334
$ws.Range("N10").Value = 0.247237
$ws.Range("O10").Value = 0.04727005612861496
$ws.Range("P10").Value = 0.04727005612861496
$ws.Range("Q10").Value = 0.02051080899077778
$ws.Range("R10").Value = 0.184597280917
$ws.Range("S10").Value = 0.0001463871063782965
$ws.Range("T10").Value = 0.0001463871063782965

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Vwf"
$ws.Range("C11").Value = "Tnfrsf11b"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.2488803333333334
$ws.Range("H11").Value = 0.7466410000000001
$ws.Range("I11").Value = 0.003096825313259592
$ws.Range("J11").Value = 0.003096825313259592
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1.661024
$ws.Range("N11").Value = 4.983072
$ws.Range("O11").Value = 0.9527299438713851
$ws.Range("P11").Value = 0.952729943871385
$ws.Range("Q11").Value = 0.4133962067946668
$ws.Range("R11").Value = 3.720565861152
$ws.Range("S11").Value = 0.002950438206881295
$ws.Range("T11").Value = 0.002950438206881295
